# fix: calander as calender
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "calander" / "Calander" typo wherever it appears on the sheet
$ws.Range("B3").Value = "calender"
$ws.Range("C3").Value = "Calender"
$ws.Range("B8").Value = "calender"
$ws.Range("C8").Value = "Calender"

# B20 previously had no value; it now holds an (empty) text value.
# A bare "" clears a cell instead of leaving an empty string behind, so
# force text entry with a leading quote (Excel's "treat as text" prefix)
# and then restore the default style so no visible formatting sticks.
$ws.Range("B20").Value = "'"
$ws.Range("B20").Style = "Normal"
